$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 859-860 (shifts the existing data, previously at
# rows 859-900, down by 2 rows to 861-902)
$ws.Rows("859:860").Insert()

# Fill in the two newly inserted rows with the new daily entries
$ws.Range("A859").Value = "'2026/02/23"
$ws.Range("B859").Value = "月"
$ws.Range("C859").Value = 23
$ws.Range("D859").Value = 201

$ws.Range("A860").Value = "'2026/02/24"
$ws.Range("B860").Value = "火"
$ws.Range("C860").Value = 2
$ws.Range("D860").Value = 201

# Re-apply the plain (un-styled) formatting of the surrounding data rows so
# the new rows match them exactly (no quote-prefix / text number format
# left behind by the literal assignments above).
$ws.Range("A858:D858").Copy()
$ws.Range("A859:D860").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
